# Auto-generated edit script applying reworked duration calculations
# and downstream Mann-Kendall trend statistics (mk_duration, mk_intra_annual),
# plus a handful of rounding updates and HCDN_2009 flag corrections on site_metrics.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("site_metrics")
$ws.Range("O2").Value = 1.432221441315653
$ws.Range("O6").Value = 0.3069872705133974
$ws.Range("AK17").Value = $true
$ws.Range("AK20").Value = $true
$ws.Range("O21").Value = 0.06371983478406092
$ws.Range("AK21").Value = $true
$ws.Range("AK23").Value = $true
$ws.Range("O29").Value = 0.01005100646505231
$ws.Range("O33").Value = 0.3094509165460451
$ws.Range("O34").Value = 0.05117192000599945
$ws.Range("O42").Value = 0.1661267140294488
$ws.Range("AK42").Value = $true
$ws.Range("N44").Value = 10.36419447219447
$ws.Range("O44").Value = 0.05048749535702617
$ws.Range("Q44").Value = 7.025
$ws.Range("O55").Value = 0.3067439527501406
$ws.Range("O57").Value = 0.3564094171642607
$ws.Range("AK58").Value = $true
$ws.Range("AK62").Value = $true
$ws.Range("N65").Value = 12.5448275862069
$ws.Range("O65").Value = 0.01993759407123781
$ws.Range("Q65").Value = 3.379310344827586
$ws.Range("O66").Value = 0.06227864556659247
$ws.Range("O71").Value = 0.01519217395941376
$ws.Range("N73").Value = 40.66
$ws.Range("O73").Value = 0.01269969739382551
$ws.Range("Q73").Value = 2.0
$ws.Range("O78").Value = 0.006459603132374887
$ws.Range("AK78").Value = $true
$ws.Range("O86").Value = 0.0518494331441559
$ws.Range("AK88").Value = $true
$ws.Range("AK91").Value = $true
$ws.Range("N95").Value = 39.04761904761905
$ws.Range("O95").Value = 0.2494472116874198
$ws.Range("Q95").Value = 2.047619047619047
$ws.Range("AK96").Value = $true
$ws.Range("AK99").Value = $true
$ws.Range("O103").Value = 0.08913716162184795
$ws.Range("O106").Value = 0.006326998972559076
$ws.Range("O111").Value = 1.123048147155037
$ws.Range("N113").Value = 8.698105158730158
$ws.Range("AK119").Value = $true
$ws.Range("O125").Value = 0.2600513475516693
$ws.Range("AK126").Value = $true
$ws.Range("O128").Value = 0.05048928385484487
$ws.Range("AK128").Value = $true
$ws.Range("AK132").Value = $true
$ws.Range("AK133").Value = $true
$ws.Range("O134").Value = 0.007966970481402706
$ws.Range("AK135").Value = $true
$ws.Range("AK136").Value = $true

$ws = $wb.Worksheets.Item("mk_duration")
$ws.Range("M4").Value = 0.9861976485245498
$ws.Range("N4").Value = 0.01729954507840292
$ws.Range("O4").Value = 0.002551020408163265
$ws.Range("P4").Value = 3.0
$ws.Range("Q4").Value = 13365.66666666667
$ws.Range("S4").Value = 13.0
$ws.Range("M18").Value = 0.6490437529262907
$ws.Range("N18").Value = -0.4550910273091032
$ws.Range("O18").Value = -0.08187134502923976
$ws.Range("P18").Value = -14.0
$ws.Range("Q18").Value = 816.0
$ws.Range("R18").Value = -0.140625
$ws.Range("S18").Value = 8.865625
$ws.Range("K22").Value = "no trend"
$ws.Range("L22").Value = $false
$ws.Range("M22").Value = 0.989125708265659
$ws.Range("N22").Value = -0.01362932551272764
$ws.Range("O22").Value = -0.003174603174603175
$ws.Range("P22").Value = -2.0
$ws.Range("Q22").Value = 5383.333333333333
$ws.Range("R22").Value = 0.0
$ws.Range("S22").Value = 6.928571428571429
$ws.Range("M26").Value = 1.0
$ws.Range("N26").Value = 0.0
$ws.Range("O26").Value = 0.0
$ws.Range("P26").Value = 0.0
$ws.Range("Q26").Value = 12865.33333333333
$ws.Range("S26").Value = 10.33333333333333
$ws.Range("M29").Value = 0.02068547487688188
$ws.Range("N29").Value = 2.313676482007082
$ws.Range("O29").Value = 0.2988505747126437
$ws.Range("P29").Value = 130.0
$ws.Range("Q29").Value = 3108.666666666667
$ws.Range("R29").Value = 0.5
$ws.Range("S29").Value = -2.055555555555555
$ws.Range("K38").Value = "no trend"
$ws.Range("L38").Value = $false
$ws.Range("M38").Value = 1.0
$ws.Range("N38").Value = 0.0
$ws.Range("O38").Value = 0.003333333333333334
$ws.Range("P38").Value = 1.0
$ws.Range("Q38").Value = 1823.666666666667
$ws.Range("S38").Value = 9.333333333333334
$ws.Range("M40").Value = 0.9196171294399185
$ws.Range("N40").Value = -0.1009160155527349
$ws.Range("O40").Value = -0.01110083256244218
$ws.Range("P40").Value = -12.0
$ws.Range("Q40").Value = 11881.33333333333
$ws.Range("R40").Value = -0.006000000000000014
$ws.Range("S40").Value = 11.638
$ws.Range("M43").Value = 0.5576025867964323
$ws.Range("N43").Value = -0.586406188193654
$ws.Range("O43").Value = -0.06312292358803986
$ws.Range("P43").Value = -57.0
$ws.Range("Q43").Value = 9119.666666666666
$ws.Range("R43").Value = -0.03333333333333333
$ws.Range("S43").Value = 8.95
$ws.Range("D44").Value = 0.0005164876856496026
$ws.Range("E44").Value = 3.472056467677751
$ws.Range("F44").Value = 0.388663967611336
$ws.Range("G44").Value = 288.0
$ws.Range("M44").Value = 0.7619667354563611
$ws.Range("N44").Value = 0.3028991285987124
$ws.Range("O44").Value = 0.03020408163265306
$ws.Range("P44").Value = 37.0
$ws.Range("M45").Value = 0.1624286847379568
$ws.Range("N45").Value = 1.396949742642015
$ws.Range("O45").Value = 0.1386054421768707
$ws.Range("P45").Value = 163.0
$ws.Range("Q45").Value = 13448.33333333333
$ws.Range("R45").Value = 0.140625
$ws.Range("S45").Value = 6.625
$ws.Range("K46").Value = "no trend"
$ws.Range("L46").Value = $false
$ws.Range("M46").Value = 0.9432068918130148
$ws.Range("N46").Value = 0.07123981788618183
$ws.Range("O46").Value = 0.01008403361344538
$ws.Range("P46").Value = 6.0
$ws.Range("Q46").Value = 4926.0
$ws.Range("R46").Value = 0.0
$ws.Range("S46").Value = 25.5
$ws.Range("K59").Value = "no trend"
$ws.Range("L59").Value = $false
$ws.Range("M59").Value = 0.2514956261995338
$ws.Range("N59").Value = 1.146724198217075
$ws.Range("O59").Value = 0.1420454545454546
$ws.Range("P59").Value = 75.0
$ws.Range("Q59").Value = 4164.333333333333
$ws.Range("R59").Value = 0.3564814814814814
$ws.Range("S59").Value = 10.96296296296297
$ws.Range("D65").Value = 0.2441369378490839
$ws.Range("E65").Value = -1.164708671158552
$ws.Range("F65").Value = -0.1551724137931035
$ws.Range("G65").Value = -63.0
$ws.Range("I65").Value = -0.2012138188608777
$ws.Range("J65").Value = 12.48366013071895
$ws.Range("M65").Value = 0.1478846168671544
$ws.Range("N65").Value = -1.447043942046255
$ws.Range("O65").Value = -0.1885057471264368
$ws.Range("P65").Value = -82.0
$ws.Range("R65").Value = -0.2333333333333334
$ws.Range("S65").Value = 13.05
$ws.Range("M66").Value = 0.9641143404717236
$ws.Range("N66").Value = 0.04499117837596322
$ws.Range("O66").Value = 0.006097560975609756
$ws.Range("P66").Value = 5.0
$ws.Range("Q66").Value = 7904.333333333333
$ws.Range("R66").Value = 0.0
$ws.Range("S66").Value = 8.25
$ws.Range("K67").Value = "no trend"
$ws.Range("L67").Value = $false
$ws.Range("M67").Value = 0.4416847963526576
$ws.Range("N67").Value = 0.7693512906668215
$ws.Range("O67").Value = 0.103448275862069
$ws.Range("P67").Value = 42.0
$ws.Range("Q67").Value = 2840.0
$ws.Range("R67").Value = 0.1083333333333333
$ws.Range("S67").Value = 6.283333333333334
$ws.Range("K70").Value = "no trend"
$ws.Range("L70").Value = $false
$ws.Range("M70").Value = 0.4592944579074156
$ws.Range("N70").Value = 0.7400091243184136
$ws.Range("O70").Value = 0.1146245059288538
$ws.Range("P70").Value = 29.0
$ws.Range("Q70").Value = 1431.666666666667
$ws.Range("R70").Value = 0.25
$ws.Range("S70").Value = 4.25
$ws.Range("K72").Value = "no trend"
$ws.Range("L72").Value = $false
$ws.Range("M72").Value = 0.5767974704841166
$ws.Range("N72").Value = 0.5580687036253706
$ws.Range("O72").Value = 0.08695652173913043
$ws.Range("P72").Value = 22.0
$ws.Range("Q72").Value = 1416.0
$ws.Range("R72").Value = 0.2
$ws.Range("S72").Value = 2.8
$ws.Range("I73").Value = -0.5931372549019608
$ws.Range("J73").Value = 32.82107843137255
$ws.Range("I95").Value = -0.4166666666666661
$ws.Range("J95").Value = 30.66666666666666
$ws.Range("M103").Value = 0.8753088059117513
$ws.Range("N103").Value = -0.1569188472540014
$ws.Range("O103").Value = -0.02016806722689076
$ws.Range("P103").Value = -12.0
$ws.Range("Q103").Value = 4914.0
$ws.Range("S103").Value = 8.0
$ws.Range("K106").Value = "no trend"
$ws.Range("L106").Value = $false
$ws.Range("M106").Value = 0.8947362490929238
$ws.Range("N106").Value = 0.1323136020412692
$ws.Range("O106").Value = 0.02153846153846154
$ws.Range("P106").Value = 7.0
$ws.Range("Q106").Value = 2056.333333333333
$ws.Range("R106").Value = 0.02272727272727273
$ws.Range("S106").Value = 7.415909090909091
$ws.Range("M113").Value = 0.112606871578502
$ws.Range("N113").Value = 1.586584096977019
$ws.Range("O113").Value = 0.1646464646464647
$ws.Range("P113").Value = 163.0
$ws.Range("Q113").Value = 10425.66666666667
$ws.Range("R113").Value = 0.0959084084084084
$ws.Range("S113").Value = 2.890015015015015
$ws.Range("M115").Value = 0.3951059323785659
$ws.Range("N115").Value = 0.8503942491228869
$ws.Range("O115").Value = 0.1096774193548387
$ws.Range("P115").Value = 51.0
$ws.Range("Q115").Value = 3457.0
$ws.Range("R115").Value = 0.05555555555555555
$ws.Range("S115").Value = 3.388888888888889
$ws.Range("K118").Value = "no trend"
$ws.Range("L118").Value = $false
$ws.Range("M118").Value = 0.7727046455311284
$ws.Range("N118").Value = 0.2888389219183989
$ws.Range("O118").Value = 0.05882352941176471
$ws.Range("P118").Value = 8.0
$ws.Range("Q118").Value = 587.3333333333334
$ws.Range("R118").Value = 0.1753246753246753
$ws.Range("S118").Value = 23.0974025974026
$ws.Range("K120").Value = "no trend"
$ws.Range("L120").Value = $false
$ws.Range("M120").Value = 0.2681558505823991
$ws.Range("N120").Value = -1.107319420872094
$ws.Range("O120").Value = -0.1448275862068966
$ws.Range("P120").Value = -63.0
$ws.Range("Q120").Value = 3135.0
$ws.Range("R120").Value = -0.2156862745098039
$ws.Range("S120").Value = 14.27745098039216
$ws.Range("M122").Value = 0.54336854876055
$ws.Range("N122").Value = -0.6077270056935287
$ws.Range("O122").Value = -0.06533776301218161
$ws.Range("P122").Value = -59.0
$ws.Range("Q122").Value = 9108.333333333334
$ws.Range("R122").Value = -0.14
$ws.Range("S122").Value = 17.94
$ws.Range("M125").Value = 0.910341484373421
$ws.Range("N125").Value = -0.1126078208221815
$ws.Range("O125").Value = -0.01724137931034483
$ws.Range("P125").Value = -7.0
$ws.Range("Q125").Value = 2839.0
$ws.Range("R125").Value = -0.03860028860028859
$ws.Range("S125").Value = 20.04040404040404
$ws.Range("M126").Value = 0.4154117810055911
$ws.Range("N126").Value = -0.8144070852755617
$ws.Range("O126").Value = -0.08405797101449275
$ws.Range("P126").Value = -87.0
$ws.Range("Q126").Value = 11151.0
$ws.Range("R126").Value = -0.1083333333333333
$ws.Range("S126").Value = 15.10416666666667
$ws.Range("M127").Value = 0.3258889674157843
$ws.Range("N127").Value = 0.9824281427292695
$ws.Range("O127").Value = 0.1433333333333333
$ws.Range("P127").Value = 43.0
$ws.Range("Q127").Value = 1827.666666666667
$ws.Range("R127").Value = 0.5238095238095238
$ws.Range("S127").Value = 20.21428571428572
$ws.Range("M138").Value = 0.3987150497726255
$ws.Range("N138").Value = 0.8439183206889768
$ws.Range("O138").Value = 0.08603145235892692
$ws.Range("P138").Value = 93.0
$ws.Range("Q138").Value = 11884.33333333333
$ws.Range("R138").Value = 0.01576923076923078
$ws.Range("S138").Value = 4.637307692307692
$ws.Range("K142").Value = "no trend"
$ws.Range("L142").Value = $false
$ws.Range("M142").Value = 0.7834235533977572
$ws.Range("N142").Value = -0.2748603200228616
$ws.Range("O142").Value = -0.03303303303303303
$ws.Range("P142").Value = -22.0
$ws.Range("Q142").Value = 5837.333333333333
$ws.Range("R142").Value = -0.03819444444444445
$ws.Range("S142").Value = 8.8875

$ws = $wb.Worksheets.Item("mk_intra_annual")
$ws.Range("M4").Value = 0.9367122131047585
$ws.Range("N4").Value = -0.07940283582325927
$ws.Range("O4").Value = -0.008503401360544218
$ws.Range("P4").Value = -10.0
$ws.Range("Q4").Value = 12847.33333333333
$ws.Range("M18").Value = 0.721946073974117
$ws.Range("N18").Value = -0.3558591171665809
$ws.Range("O18").Value = -0.06432748538011696
$ws.Range("P18").Value = -11.0
$ws.Range("Q18").Value = 789.6666666666666
$ws.Range("S18").Value = 4.0
$ws.Range("K22").Value = "no trend"
$ws.Range("L22").Value = $false
$ws.Range("M22").Value = 0.9124127597319072
$ws.Range("N22").Value = 0.1099957320192416
$ws.Range("O22").Value = 0.01428571428571429
$ws.Range("P22").Value = 9.0
$ws.Range("Q22").Value = 5289.666666666667
$ws.Range("R22").Value = 0.0
$ws.Range("S22").Value = 4.5
$ws.Range("M26").Value = 0.5438898774467631
$ws.Range("N26").Value = -0.6069412864988872
$ws.Range("O26").Value = -0.0586734693877551
$ws.Range("P26").Value = -69.0
$ws.Range("Q26").Value = 12552.33333333333
$ws.Range("M29").Value = 0.6760665590790988
$ws.Range("N29").Value = -0.4178366412198368
$ws.Range("O29").Value = -0.05517241379310345
$ws.Range("P29").Value = -24.0
$ws.Range("Q29").Value = 3030.0
$ws.Range("K38").Value = "no trend"
$ws.Range("L38").Value = $false
$ws.Range("M38").Value = 0.567965060874607
$ws.Range("N38").Value = 0.5710510169959229
$ws.Range("O38").Value = 0.08333333333333333
$ws.Range("P38").Value = 25.0
$ws.Range("Q38").Value = 1766.333333333333
$ws.Range("S38").Value = 3.0
$ws.Range("M40").Value = 0.5094322249010004
$ws.Range("N40").Value = -0.6597220328454717
$ws.Range("O40").Value = -0.06567992599444958
$ws.Range("P40").Value = -71.0
$ws.Range("Q40").Value = 11258.33333333333
$ws.Range("S40").Value = 3.0
$ws.Range("M43").Value = 0.7487382613128855
$ws.Range("N43").Value = 0.3203035117238672
$ws.Range("O43").Value = 0.03322259136212625
$ws.Range("P43").Value = 30.0
$ws.Range("Q43").Value = 8197.333333333334
$ws.Range("S43").Value = 3.0
$ws.Range("I44").Value = -0.3333333333333333
$ws.Range("J44").Value = 11.33333333333333
$ws.Range("R44").Value = -0.2692307692307692
$ws.Range("S44").Value = 10.59615384615385
$ws.Range("M45").Value = 0.1099969748728487
$ws.Range("N45").Value = -1.598206737151498
$ws.Range("O45").Value = -0.1547619047619048
$ws.Range("P45").Value = -182.0
$ws.Range("Q45").Value = 12826.0
$ws.Range("K46").Value = "no trend"
$ws.Range("L46").Value = $false
$ws.Range("M46").Value = 0.3477560684617595
$ws.Range("N46").Value = -0.9389506772171409
$ws.Range("O46").Value = -0.09915966386554621
$ws.Range("P46").Value = -59.0
$ws.Range("Q46").Value = 3815.666666666667
$ws.Range("R46").Value = 0.0
$ws.Range("S46").Value = 1.0
$ws.Range("K59").Value = "no trend"
$ws.Range("L59").Value = $false
$ws.Range("M59").Value = 0.6372999986566994
$ws.Range("N59").Value = 0.4714771987049622
$ws.Range("O59").Value = 0.05681818181818182
$ws.Range("P59").Value = 30.0
$ws.Range("Q59").Value = 3783.333333333333
$ws.Range("R59").Value = 0.0
$ws.Range("S59").Value = 2.0
$ws.Range("D65").Value = 0.07954730198972015
$ws.Range("E65").Value = 1.753318789833714
$ws.Range("F65").Value = 0.2216748768472906
$ws.Range("G65").Value = 90.0
$ws.Range("H65").Value = 2576.666666666667
$ws.Range("M65").Value = 0.1918233548641619
$ws.Range("N65").Value = 1.305204110593214
$ws.Range("O65").Value = 0.1632183908045977
$ws.Range("P65").Value = 71.0
$ws.Range("Q65").Value = 2876.333333333333
$ws.Range("K66").Value = "no trend"
$ws.Range("L66").Value = $false
$ws.Range("M66").Value = 0.3444639880539495
$ws.Range("N66").Value = 0.9453818046781428
$ws.Range("O66").Value = 0.1024390243902439
$ws.Range("P66").Value = 84.0
$ws.Range("Q66").Value = 7708.0
$ws.Range("R66").Value = 0.0
$ws.Range("S66").Value = 4.0
$ws.Range("K67").Value = "no trend"
$ws.Range("L67").Value = $false
$ws.Range("M67").Value = 0.8788010010727179
$ws.Range("N67").Value = -0.1524893355507309
$ws.Range("O67").Value = -0.02216748768472906
$ws.Range("P67").Value = -9.0
$ws.Range("Q67").Value = 2752.333333333333
$ws.Range("R67").Value = 0.0
$ws.Range("S67").Value = 4.0
$ws.Range("K70").Value = "no trend"
$ws.Range("L70").Value = $false
$ws.Range("M70").Value = 0.0797981882080041
$ws.Range("N70").Value = -1.751858231685006
$ws.Range("O70").Value = -0.2608695652173913
$ws.Range("P70").Value = -66.0
$ws.Range("Q70").Value = 1376.666666666667
$ws.Range("R70").Value = -0.1
$ws.Range("S70").Value = 4.1
$ws.Range("K72").Value = "no trend"
$ws.Range("L72").Value = $false
$ws.Range("M72").Value = 0.3455232426259909
$ws.Range("N72").Value = -0.9433082785996189
$ws.Range("O72").Value = -0.1422924901185771
$ws.Range("P72").Value = -36.0
$ws.Range("Q72").Value = 1376.666666666667
$ws.Range("R72").Value = -0.05882352941176471
$ws.Range("S72").Value = 2.647058823529412
$ws.Range("D73").Value = 0.02267985566984976
$ws.Range("E73").Value = -2.278784887918315
$ws.Range("F73").Value = -0.3115942028985507
$ws.Range("G73").Value = -86.0
$ws.Range("I73").Value = -0.05409356725146199
$ws.Range("J73").Value = 2.122076023391813
$ws.Range("M73").Value = 0.00007722878086746654
$ws.Range("N73").Value = -3.95284092599084
$ws.Range("O73").Value = -0.3542857142857143
$ws.Range("P73").Value = -434.0
$ws.Range("R73").Value = -0.02564102564102564
$ws.Range("S73").Value = 0.6282051282051282
$ws.Range("D95").Value = 0.2170540149270519
$ws.Range("E95").Value = -1.234399661981982
$ws.Range("F95").Value = -0.1857142857142857
$ws.Range("G95").Value = -39.0
$ws.Range("H95").Value = 947.6666666666666
$ws.Range("M95").Value = 0.756441108307621
$ws.Range("N95").Value = -0.3101576047729196
$ws.Range("O95").Value = -0.03373819163292847
$ws.Range("P95").Value = -25.0
$ws.Range("Q95").Value = 5987.666666666667
$ws.Range("M103").Value = 0.5262476380508285
$ws.Range("N103").Value = -0.6337444087790188
$ws.Range("O103").Value = -0.07563025210084033
$ws.Range("P103").Value = -45.0
$ws.Range("Q103").Value = 4820.333333333333
$ws.Range("S103").Value = 3.0
$ws.Range("K106").Value = "no trend"
$ws.Range("L106").Value = $false
$ws.Range("M106").Value = 0.720045021432083
$ws.Range("N106").Value = -0.3583986237519272
$ws.Range("O106").Value = -0.05230769230769231
$ws.Range("P106").Value = -17.0
$ws.Range("Q106").Value = 1993.0
$ws.Range("S106").Value = 3.0
$ws.Range("K113").Value = "decreasing"
$ws.Range("L113").Value = $true
$ws.Range("M113").Value = 0.04948737124695035
$ws.Range("N113").Value = -1.964368515926561
$ws.Range("O113").Value = -0.202020202020202
$ws.Range("P113").Value = -200.0
$ws.Range("Q113").Value = 10262.66666666667
$ws.Range("R113").Value = -0.07362240289069558
$ws.Range("S113").Value = 5.619692863595303
$ws.Range("M115").Value = 0.3282060661285529
$ws.Range("N115").Value = -0.9777336665949743
$ws.Range("O115").Value = -0.1247311827956989
$ws.Range("P115").Value = -58.0
$ws.Range("Q115").Value = 3398.666666666667
$ws.Range("R115").Value = -0.05263157894736842
$ws.Range("S115").Value = 6.789473684210527
$ws.Range("K118").Value = "no trend"
$ws.Range("L118").Value = $false
$ws.Range("M118").Value = 1.0
$ws.Range("N118").Value = 0.0
$ws.Range("O118").Value = -0.007352941176470588
$ws.Range("P118").Value = -1.0
$ws.Range("Q118").Value = 515.6666666666666
$ws.Range("S118").Value = 2.0
$ws.Range("K120").Value = "no trend"
$ws.Range("L120").Value = $false
$ws.Range("M120").Value = 0.1547212957319766
$ws.Range("N120").Value = 1.42305126113676
$ws.Range("O120").Value = 0.1816091954022989
$ws.Range("P120").Value = 79.0
$ws.Range("Q120").Value = 3004.333333333333
$ws.Range("R120").Value = 0.04347826086956522
$ws.Range("S120").Value = 2.369565217391304
$ws.Range("M122").Value = 0.6090473046930951
$ws.Range("N122").Value = 0.5114338394653799
$ws.Range("O122").Value = 0.05315614617940199
$ws.Range("P122").Value = 48.0
$ws.Range("Q122").Value = 8445.333333333334
$ws.Range("M125").Value = 0.4214420309678952
$ws.Range("N125").Value = 0.8039219871759481
$ws.Range("O125").Value = 0.1009852216748768
$ws.Range("P125").Value = 41.0
$ws.Range("Q125").Value = 2475.666666666667
$ws.Range("S125").Value = 2.0
$ws.Range("K126").Value = "no trend"
$ws.Range("L126").Value = $false
$ws.Range("M126").Value = 0.06818489781853465
$ws.Range("N126").Value = 1.823782914899107
$ws.Range("O126").Value = 0.1835748792270532
$ws.Range("P126").Value = 190.0
$ws.Range("Q126").Value = 10739.33333333333
$ws.Range("R126").Value = 0.0303030303030303
$ws.Range("S126").Value = 2.318181818181818
$ws.Range("M127").Value = 0.09365377837159627
$ws.Range("N127").Value = 1.676431115802568
$ws.Range("O127").Value = 0.22
$ws.Range("P127").Value = 66.0
$ws.Range("Q127").Value = 1503.333333333333
$ws.Range("S127").Value = 2.0
$ws.Range("M138").Value = 0.8968316034995198
$ws.Range("N138").Value = -0.1296648371974174
$ws.Range("O138").Value = -0.01387604070305273
$ws.Range("P138").Value = -15.0
$ws.Range("Q138").Value = 11657.66666666667
$ws.Range("K142").Value = "no trend"
$ws.Range("L142").Value = $false
$ws.Range("M142").Value = 0.124767437701691
$ws.Range("N142").Value = 1.535066729490695
$ws.Range("O142").Value = 0.1756756756756757
$ws.Range("P142").Value = 117.0
$ws.Range("Q142").Value = 5710.333333333333
$ws.Range("R142").Value = 0.04880952380952381
$ws.Range("S142").Value = 3.121428571428571

